# Auto-generated edit script: updates Leve profit-tracking values across sheets
# (scheduled-runner style refresh of currentAveragePrice / LevePrice / LeveProfit columns)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 250
$ws.Range("I10").Value = 250
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 250
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = 43
$ws.Range("N10").ClearContents()
$ws.Range("H13").Value = 1013.1667
$ws.Range("J13").Value = 1086.3334
$ws.Range("L13").Value = 1086.3334
$ws.Range("N13").Value = -1424.3334
$ws.Range("H116").Value = 1995.6666
$ws.Range("I116").Value = 1993.5
$ws.Range("K116").Value = 1993.5
$ws.Range("M116").Value = 1448.5
$ws.Range("H131").Value = 1390.4
$ws.Range("I131").Value = 1390.4
$ws.Range("K131").Value = 4171.200000000001
$ws.Range("M131").Value = 868.7999999999993
$ws.Range("H137").Value = 3434.3845
$ws.Range("I137").Value = 2922.6365
$ws.Range("K137").Value = 8767.9095
$ws.Range("M137").Value = -6217.9095

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3895.6
$ws.Range("I2").Value = 3620
$ws.Range("J2").Value = 4998
$ws.Range("K2").Value = 3620
$ws.Range("L2").Value = 4998
$ws.Range("M2").Value = -3507
$ws.Range("N2").Value = -5224
$ws.Range("H32").Value = 1945.931
$ws.Range("I32").Value = 1982.7407
$ws.Range("K32").Value = 1982.7407
$ws.Range("M32").Value = -1695.7407
$ws.Range("H74").Value = 3753.0667
$ws.Range("J74").Value = 5899.8335
$ws.Range("L74").Value = 5899.8335
$ws.Range("N74").Value = -7647.8335
$ws.Range("H77").Value = 3753.0667
$ws.Range("J77").Value = 5899.8335
$ws.Range("L77").Value = 29499.1675
$ws.Range("N77").Value = -38235.1675
$ws.Range("H116").Value = 3895.6
$ws.Range("I116").Value = 3620
$ws.Range("J116").Value = 4998
$ws.Range("K116").Value = 3620
$ws.Range("L116").Value = 4998
$ws.Range("M116").Value = -1326
$ws.Range("N116").Value = -9586
$ws.Range("H122").Value = 2971.2856
$ws.Range("I122").Value = 3199.75
$ws.Range("J122").Value = 2666.6667
$ws.Range("K122").Value = 9599.25
$ws.Range("L122").Value = 8000.000100000001
$ws.Range("M122").Value = -7149.25
$ws.Range("N122").Value = -12900.0001
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3895.6
$ws.Range("I3").Value = 3620
$ws.Range("J3").Value = 4998
$ws.Range("K3").Value = 3620
$ws.Range("L3").Value = 4998
$ws.Range("M3").Value = -3506
$ws.Range("N3").Value = -5226
$ws.Range("H105").Value = 3354.9167
$ws.Range("I105").Value = 2608.4285
$ws.Range("J105").Value = 4400
$ws.Range("K105").Value = 2608.4285
$ws.Range("L105").Value = 4400
$ws.Range("M105").Value = -861.4285
$ws.Range("N105").Value = -7894
$ws.Range("H122").Value = 100000
$ws.Range("J122").Value = 100000
$ws.Range("L122").Value = 100000
$ws.Range("N122").Value = -109800

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 276
$ws.Range("I2").Value = 351.33334
$ws.Range("J2").Value = 50
$ws.Range("K2").Value = 351.33334
$ws.Range("L2").Value = 50
$ws.Range("M2").Value = -238.33334
$ws.Range("N2").Value = -276
$ws.Range("H5").Value = 293.625
$ws.Range("I5").Value = 76.666664
$ws.Range("J5").Value = 423.8
$ws.Range("K5").Value = 76.666664
$ws.Range("L5").Value = 423.8
$ws.Range("M5").Value = 35.333336
$ws.Range("N5").Value = -647.8
$ws.Range("H10").Value = 303.42856
$ws.Range("I10").Value = 319.83334
$ws.Range("J10").Value = 205
$ws.Range("K10").Value = 319.83334
$ws.Range("L10").Value = 205
$ws.Range("M10").Value = -180.83334
$ws.Range("N10").Value = -483
$ws.Range("H11").Value = 205
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 205
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 205
$ws.Range("M11").ClearContents()
$ws.Range("N11").Value = -485
$ws.Range("H31").Value = 6340.65
$ws.Range("I31").Value = 4618.1113
$ws.Range("K31").Value = 4618.1113
$ws.Range("M31").Value = -4323.1113
$ws.Range("H34").Value = 6340.65
$ws.Range("I34").Value = 4618.1113
$ws.Range("K34").Value = 4618.1113
$ws.Range("M34").Value = -4416.1113
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()
$ws.Range("H86").Value = 4300
$ws.Range("I86").Value = 4300
$ws.Range("K86").Value = 4300
$ws.Range("M86").Value = -3177
$ws.Range("H89").Value = 4300
$ws.Range("I89").Value = 4300
$ws.Range("K89").Value = 21500
$ws.Range("M89").Value = -15884

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 7196.625
$ws.Range("I102").Value = 7082.5713
$ws.Range("J102").Value = 7995
$ws.Range("K102").Value = 7082.5713
$ws.Range("L102").Value = 7995
$ws.Range("M102").Value = -5460.5713
$ws.Range("N102").Value = -11239
$ws.Range("H126").Value = 5198
$ws.Range("I126").Value = 5254.5713
$ws.Range("K126").Value = 15763.7139
$ws.Range("M126").Value = -13293.7139

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4073
$ws.Range("I40").Value = 3751.8333
$ws.Range("K40").Value = 3751.8333
$ws.Range("M40").Value = -3615.8333
$ws.Range("H46").Value = 1571.5714
$ws.Range("I46").Value = 1733
$ws.Range("J46").Value = 1450.5
$ws.Range("K46").Value = 1733
$ws.Range("L46").Value = 1450.5
$ws.Range("M46").Value = -1545
$ws.Range("N46").Value = -1826.5
$ws.Range("H61").Value = 3563.25
$ws.Range("I61").Value = 4082.6667
$ws.Range("J61").Value = 2005
$ws.Range("K61").Value = 4082.6667
$ws.Range("L61").Value = 2005
$ws.Range("M61").Value = -3880.6667
$ws.Range("N61").Value = -2409
$ws.Range("H68").Value = 2074.6365
$ws.Range("I68").Value = 2157.1428
$ws.Range("J68").Value = 1930.25
$ws.Range("K68").Value = 2157.1428
$ws.Range("L68").Value = 1930.25
$ws.Range("M68").Value = -1408.1428
$ws.Range("N68").Value = -3428.25
$ws.Range("H71").Value = 2074.6365
$ws.Range("I71").Value = 2157.1428
$ws.Range("J71").Value = 1930.25
$ws.Range("K71").Value = 10785.714
$ws.Range("L71").Value = 9651.25
$ws.Range("M71").Value = -7041.714
$ws.Range("N71").Value = -17139.25
$ws.Range("H82").Value = 2368.75
$ws.Range("I82").Value = 1991.6666
$ws.Range("K82").Value = 1991.6666
$ws.Range("M82").Value = -1630.6666
$ws.Range("H85").Value = 2368.75
$ws.Range("I85").Value = 1991.6666
$ws.Range("K85").Value = 1991.6666
$ws.Range("M85").Value = -743.6666
$ws.Range("H113").Value = 3563.25
$ws.Range("I113").Value = 4082.6667
$ws.Range("J113").Value = 2005
$ws.Range("K113").Value = 4082.6667
$ws.Range("L113").Value = 2005
$ws.Range("M113").Value = -1912.6667
$ws.Range("N113").Value = -6345
$ws.Range("H118").Value = 99995
$ws.Range("J118").Value = 99995
$ws.Range("L118").Value = 99995
$ws.Range("N118").Value = -103309

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1037
$ws.Range("I113").Value = 448.2
$ws.Range("J113").Value = 2018.3334
$ws.Range("K113").Value = 1344.6
$ws.Range("L113").Value = 6055.0002
$ws.Range("M113").Value = 825.4000000000001
$ws.Range("N113").Value = -10395.0002
